$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column F: "BF01" header + =1/E{row} values ------------------------
# Header cell F1 picks up the same header look (font/border/valign) as the
# other header cells by copying A1's format, then setting its text.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F1").Value = "BF01"

# F2 gets its own (non-shared) formula; F3:F15 filled in one shot so the
# engine groups them into a shared formula, matching the original edit.
$ws.Range("F2").Formula = "=1/E2"
$ws.Range("F3:F15").Formula = "=1/E3"

# --- Header emphasis: B1:D1 become italic -----------------------------------
$ws.Range("B1:D1").Font.Italic = $true

# --- Row 1 height tweak ------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15.75

# --- Selection / view state to match the saved workbook ---------------------
$ws.Range("F6:F8").Select()
